# Auto-generated: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.875.51'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.623.15'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.995'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.90'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.79'
$ws.Range("E8").Value = '  +11.22%  '
$ws.Range("E9").Value = '  +3.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0610'
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").Value = '1.854.89'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("D13").Value = '1.619.60'
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.569'
$ws.Range("E14").Value = '  +6.00%  '
$ws.Range("E15").Value = '  +4.88%  '
$ws.Range("D16").Value = '29.907.71'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("E17").Value = '  +15.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.56'
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.13'
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  +3.43%  '
$ws.Range("E23").Value = '  +4.22%  '
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.77'
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.67'
$ws.Range("E26").Value = '  +2.69%  '
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.59'
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  +3.46%  '
$ws.Range("E31").Value = '  +5.41%  '
$ws.Range("E32").Value = '  +3.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").Value = '1.424.81'
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.65'
$ws.Range("E35").Value = '  +7.20%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("E40").Value = '  +3.35%  '
$ws.Range("E41").Value = '  +2.90%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.833'
$ws.Range("E43").Value = '  +4.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.30'
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.04'
$ws.Range("E45").Value = '  +4.89%  '
$ws.Range("E46").Value = '  +19.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.995'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.41'
$ws.Range("E48").Value = '  +2.52%  '
$ws.Range("D49").Value = '1.763.54'
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.07'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("E51").Value = '  +5.94%  '
